$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-08-30 21:01:36"

# --- Rows 4-6: the scraper re-ran and the product order shifted by one
# (new row4 = old row6 "Tena", new row5 = old row4 "Avela", new row6 = old
# row5 "Naturaline"). Rotate via Cut/Paste (through a scratch row) so the
# original cell types/values travel with the row instead of being retyped,
# which would coerce numeric-looking text (ids) into real numbers.
$ws.Range("A4:O4").Cut($ws.Range("A100:O100"))
$ws.Range("A6:O6").Cut($ws.Range("A4:O4"))
$ws.Range("A5:O5").Cut($ws.Range("A6:O6"))
$ws.Range("A100:O100").Cut($ws.Range("A5:O5"))

# --- All data rows: the timestamp column reflects the re-scrape time ---
$ws.Range("O2:O65").Value = $newTimestamp
